$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 9361.362999999999
$ws.Range("I12").Value = 11302.556
$ws.Range("K12").Value = 11302.556
$ws.Range("M12").Value = -11132.556
$ws.Range("H19").Value = 3832.5
$ws.Range("I19").Value = 5800
$ws.Range("K19").Value = 5800
$ws.Range("M19").Value = -5625
$ws.Range("H86").Value = 125000720
$ws.Range("I86").Value = 200000340
$ws.Range("J86").Value = 1348.3334
$ws.Range("K86").Value = 200000340
$ws.Range("L86").Value = 1348.3334
$ws.Range("M86").Value = -199999217
$ws.Range("N86").Value = -3594.3334
$ws.Range("H89").Value = 125000720
$ws.Range("I89").Value = 200000340
$ws.Range("J89").Value = 1348.3334
$ws.Range("K89").Value = 1000001700
$ws.Range("L89").Value = 6741.666999999999
$ws.Range("M89").Value = -999996084
$ws.Range("N89").Value = -17973.667
$ws.Range("H112").Value = 4209.724
$ws.Range("J112").Value = 4209.724
$ws.Range("L112").Value = 12629.172
$ws.Range("N112").Value = -14845.172
$ws.Range("H116").Value = 23619830
$ws.Range("I116").Value = 26990852
$ws.Range("J116").Value = 22664.666
$ws.Range("K116").Value = 26990852
$ws.Range("L116").Value = 22664.666
$ws.Range("M116").Value = -26987410
$ws.Range("N116").Value = -29548.666
$ws.Range("H132").Value = 7580.4507
$ws.Range("I132").Value = 4343.2256
$ws.Range("K132").Value = 13029.6768
$ws.Range("M132").Value = -10499.6768
$ws.Range("H137").Value = 188272.22
$ws.Range("I137").Value = 252833.67
$ws.Range("K137").Value = 758501.01
$ws.Range("M137").Value = -755951.01
$ws.Range("H138").Value = 4050.36
$ws.Range("I138").Value = 857.9149
$ws.Range("J138").Value = 6881.396
$ws.Range("K138").Value = 2573.7447
$ws.Range("L138").Value = 20644.188
$ws.Range("M138").Value = 2566.2553
$ws.Range("N138").Value = -30924.188
$ws.Range("H141").Value = 2476.7
$ws.Range("I141").Value = 2303.3584
$ws.Range("K141").Value = 6910.0752
$ws.Range("M141").Value = -1730.0752

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20563.76
$ws.Range("I32").Value = 24390.258
$ws.Range("K32").Value = 24390.258
$ws.Range("M32").Value = -24103.258
$ws.Range("H45").Value = 2351.2
$ws.Range("I45").Value = 1242.4
$ws.Range("K45").Value = 1242.4
$ws.Range("M45").Value = -865.4000000000001
$ws.Range("H122").Value = 4647.1875
$ws.Range("I122").Value = 3367.3914
$ws.Range("J122").Value = 7917.778
$ws.Range("K122").Value = 10102.1742
$ws.Range("L122").Value = 23753.334
$ws.Range("M122").Value = -7652.174199999999
$ws.Range("N122").Value = -28653.334
$ws.Range("H132").Value = 29921.863
$ws.Range("I132").Value = 53371.184
$ws.Range("J132").Value = 6472.5454
$ws.Range("K132").Value = 160113.552
$ws.Range("L132").Value = 19417.6362
$ws.Range("M132").Value = -157583.552
$ws.Range("N132").Value = -24477.6362

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1197.5
$ws.Range("I94").Value = 1292.7273
$ws.Range("K94").Value = 1292.7273
$ws.Range("M94").Value = -841.7273
$ws.Range("H134").Value = 4472.9165
$ws.Range("I134").Value = 3113.75
$ws.Range("K134").Value = 9341.25
$ws.Range("M134").Value = -6806.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2204.2415
$ws.Range("I16").Value = 1932.44
$ws.Range("J16").Value = 3903
$ws.Range("K16").Value = 1932.44
$ws.Range("L16").Value = 3903
$ws.Range("M16").Value = -1645.44
$ws.Range("N16").Value = -4477
$ws.Range("H22").Value = 676.4
$ws.Range("I22").Value = 716.36365
$ws.Range("K22").Value = 716.36365
$ws.Range("M22").Value = -366.36365
$ws.Range("H31").Value = 12199059
$ws.Range("I31").Value = 17860638
$ws.Range("J31").Value = 4889.5386
$ws.Range("K31").Value = 17860638
$ws.Range("L31").Value = 4889.5386
$ws.Range("M31").Value = -17860343
$ws.Range("N31").Value = -5479.5386
$ws.Range("H34").Value = 12199059
$ws.Range("I34").Value = 17860638
$ws.Range("J34").Value = 4889.5386
$ws.Range("K34").Value = 17860638
$ws.Range("L34").Value = 4889.5386
$ws.Range("M34").Value = -17860436
$ws.Range("N34").Value = -5293.5386
$ws.Range("H41").Value = 17587
$ws.Range("I41").Value = 3449.3333
$ws.Range("K41").Value = 3449.3333
$ws.Range("M41").Value = -3021.3333
$ws.Range("H99").Value = 13119.581
$ws.Range("I99").Value = 14834.643
$ws.Range("K99").Value = 14834.643
$ws.Range("M99").Value = -13336.643
$ws.Range("H113").Value = 2204.2415
$ws.Range("I113").Value = 1932.44
$ws.Range("J113").Value = 3903
$ws.Range("K113").Value = 1932.44
$ws.Range("L113").Value = 3903
$ws.Range("M113").Value = 237.5599999999999
$ws.Range("N113").Value = -8243
$ws.Range("H126").Value = 13119.581
$ws.Range("I126").Value = 14834.643
$ws.Range("K126").Value = 44503.929
$ws.Range("M126").Value = -42033.929
$ws.Range("H134").Value = 1401.5
$ws.Range("I134").Value = 1176.7916
$ws.Range("K134").Value = 3530.3748
$ws.Range("M134").Value = -995.3748000000001
$ws.Range("H141").Value = 81759.31
$ws.Range("J141").Value = 80897.3
$ws.Range("L141").Value = 80897.3
$ws.Range("N141").Value = -91257.3

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 974.3333
$ws.Range("J113").Value = 969.4
$ws.Range("L113").Value = 2908.2
$ws.Range("N113").Value = -7248.2
$ws.Range("H122").Value = 822.5
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 596.6667
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 5370.0003
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -10270.0003

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 84847.07000000001
$ws.Range("I80").Value = 96623.55
$ws.Range("K80").Value = 96623.55
$ws.Range("M80").Value = -95625.55
$ws.Range("H83").Value = 84847.07000000001
$ws.Range("I83").Value = 96623.55
$ws.Range("K83").Value = 483117.75
$ws.Range("M83").Value = -478125.75
$ws.Range("H126").Value = 5416.6294
$ws.Range("I126").Value = 5033
$ws.Range("K126").Value = 15099
$ws.Range("M126").Value = -12629
$ws.Range("H132").Value = 357606
$ws.Range("I132").Value = 83984.24000000001
$ws.Range("K132").Value = 251952.72
$ws.Range("M132").Value = -249422.72

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 37999
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()  # was -38590
$ws.Range("H68").Value = 3644.7334
$ws.Range("I68").Value = 2723.1667
$ws.Range("J68").Value = 7331
$ws.Range("K68").Value = 2723.1667
$ws.Range("L68").Value = 7331
$ws.Range("M68").Value = -1974.1667
$ws.Range("N68").Value = -8829
$ws.Range("H71").Value = 3644.7334
$ws.Range("I71").Value = 2723.1667
$ws.Range("J71").Value = 7331
$ws.Range("K71").Value = 13615.8335
$ws.Range("L71").Value = 36655
$ws.Range("M71").Value = -9871.833500000001
$ws.Range("N71").Value = -44143
$ws.Range("H122").Value = 4822.7
$ws.Range("I122").Value = 1998.6875
$ws.Range("K122").Value = 5996.0625
$ws.Range("M122").Value = -3546.0625
$ws.Range("H132").Value = 1730.96
$ws.Range("I132").Value = 1781.9584
$ws.Range("J132").Value = 507
$ws.Range("K132").Value = 5345.8752
$ws.Range("L132").Value = 1521
$ws.Range("M132").Value = -2815.8752
$ws.Range("N132").Value = -6581
$ws.Range("H136").Value = 2362.525
$ws.Range("I136").Value = 1508.9846
$ws.Range("J136").Value = 6061.2
$ws.Range("K136").Value = 4526.9538
$ws.Range("L136").Value = 18183.6
$ws.Range("M136").Value = -1976.9538
$ws.Range("N136").Value = -23283.6
$ws.Range("H139").Value = 81497
$ws.Range("J139").Value = 81497
$ws.Range("L139").Value = 81497
$ws.Range("N139").Value = -91777

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2981.3333
$ws.Range("J62").Value = 2995.6667
$ws.Range("L62").Value = 2995.6667
$ws.Range("N62").Value = -4243.6667
$ws.Range("H65").Value = 2981.3333
$ws.Range("J65").Value = 2995.6667
$ws.Range("L65").Value = 14978.3335
$ws.Range("N65").Value = -21218.3335
$ws.Range("H122").Value = 4023.2632
$ws.Range("I122").Value = 2453.25
$ws.Range("K122").Value = 7359.75
$ws.Range("M122").Value = -4909.75
$ws.Range("H132").Value = 1992.3433
$ws.Range("I132").Value = 836.6957
$ws.Range("K132").Value = 2510.0871
$ws.Range("M132").Value = 19.91290000000026
$ws.Range("H136").Value = 5355.42
$ws.Range("I136").Value = 1161.3877
$ws.Range("J136").Value = 9384.98
$ws.Range("K136").Value = 3484.1631
$ws.Range("L136").Value = 28154.94
$ws.Range("M136").Value = -934.1630999999998
$ws.Range("N136").Value = -33254.94
$ws.Range("H140").Value = 75441.836
$ws.Range("J140").Value = 75441.836
$ws.Range("L140").Value = 75441.836
$ws.Range("N140").Value = -85801.836
